$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (Changed) date column C for all data rows (2-482)
# from serial date 45180 (2023-09-11) to 45181 (2023-09-12).
$lastRow = 482
$ws.Range("C2:C$lastRow").Value = 45181
